$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 155, shifting existing rows 155-191 down to 156-192.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new data record.
$ws.Cells.Item(155, 1).Value = 8
$ws.Cells.Item(155, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 45093
$ws.Cells.Item(155, 5).Value = 4
$ws.Cells.Item(155, 6).Value = 100112052
$ws.Cells.Item(155, 7).Value = "Albahaca"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 800
$ws.Cells.Item(155, 11).Value = 2800
$ws.Cells.Item(155, 12).Value = 3000
$ws.Cells.Item(155, 13).Value = 2900
$ws.Cells.Item(155, 14).Value = "$/paquete"
$ws.Cells.Item(155, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value = 2900
$ws.Cells.Item(155, 17).Value = 1
$ws.Cells.Item(155, 18).Value = "Hortaliza"
